# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 633
$ws1.Range("F6").Value = 9579
$ws1.Range("F7").Value = 860
$ws1.Range("F9").Value = 1211
$ws1.Range("F10").Value = 2211
$ws1.Range("F18").Value = 1333

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 633
$ws4.Range("F7").Value = 9579
$ws4.Range("F8").Value = 860
$ws4.Range("F10").Value = 1211
$ws4.Range("F11").Value = 2211
$ws4.Range("F19").Value = 1333

$wb.Save()
